# Edit script for "25/09/2017 CHITRA AND MAMATHA CHICK IN"
#
# 1) Merge the two runs that make up "Sun Sep 23 11:19:49 PDT 2017" into a
#    single run (no visible text change).
# 2) Append a new purchase-details entry (Mon Sep 24 ...) after the last
#    existing entry in the document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: collapse "Sun Sep 23" + " 11:19:49 PDT 2017" into one run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Sun Sep 23 11:19:49 PDT 2017", $false, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Sun Sep 23 11:19:49 PDT 2017", 2) | Out-Null

# ---------------------------------------------------------------------
# Step 2: locate the anchor paragraph -- the bold "Amount balance" line
# that immediately follows "Total Price ... - 4202.0" -- and append the
# new block of paragraphs after it.
# ---------------------------------------------------------------------
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like "Amount balance*4202.0*") {
        $anchorIndex = $i
    }
}

function Add-Para($cur, [string]$text, [int]$bold) {
    $cur.Range.InsertParagraphAfter()
    $idx = $cur.Index + 1
    $n = $d.Paragraphs($idx)
    $n.Range.Bold = $bold
    if ($text.Length -gt 0) {
        $n.Range.Text = $text
    }
    return $n
}

$cur = $d.Paragraphs($anchorIndex)
$cur = Add-Para $cur "" 1
$cur = Add-Para $cur "Mon Sep 24 10:51:24 PDT 2017" 0
$cur = Add-Para $cur "Person Name`t`t`t`t- MAHADEVA" 0
$cur = Add-Para $cur "---------------------------------------------------------------" 0
$cur = Add-Para $cur "Item Name`t`t`t`t- POTATO" 0
$cur = Add-Para $cur "Number of Pockets`t`t`t- 1" 0
$cur = Add-Para $cur "Number of KGs`t`t`t- 54" 0
$cur = Add-Para $cur "Rate`t`t`t`t`t- 13" 0
$cur = Add-Para $cur "Total Price`t`t`t`t- 702.0" 0
$cur = Add-Para $cur "Amount balance`t`t`t- 4904.0" 1
$cur = Add-Para $cur "" 0
$cur = Add-Para $cur "" 1
